# Germany Verbandsliga - update odds database (commit: "Atualização de bases
# das ligas, do dia: 19-04-2024 às 23:27").
#
# The upstream source re-sorted a handful of fixtures that share the same
# kick-off date/time (their stable sort key), which swaps the full row
# payload (every column except the running "id" in column A) between a
# few adjacent rows, and appends two freshly scraped fixtures at the end
# of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: swap the B:AC payload of two rows, leaving column A (the
# sequential "id") untouched.
# ---------------------------------------------------------------------
function Swap-RowPayload {
    param($ws, [int]$rowA, [int]$rowB)

    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")

    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2

    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

# Simple adjacent-row swaps (two fixtures trade full content).
Swap-RowPayload $ws 3 4
Swap-RowPayload $ws 16 17
Swap-RowPayload $ws 20 21
Swap-RowPayload $ws 46 47
Swap-RowPayload $ws 69 70
Swap-RowPayload $ws 117 118

# Three-way rotation among rows 5, 6, 7: new5 = old7, new6 = old5, new7 = old6.
$r5 = $ws.Range("B5:AC5").Value2
$r6 = $ws.Range("B6:AC6").Value2
$r7 = $ws.Range("B7:AC7").Value2

$ws.Range("B5:AC5").Value2 = $r7
$ws.Range("B6:AC6").Value2 = $r5
$ws.Range("B7:AC7").Value2 = $r6

# ---------------------------------------------------------------------
# Append two newly scraped fixtures as rows 139 / 140.
# ---------------------------------------------------------------------
function Set-RowValues {
    param($ws, [int]$row, [object[]]$values)

    $arr = New-Object 'object[,]' 1, $values.Length
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0, $i] = $values[$i]
    }
    $ws.Range("B$row`:AC$row").Value2 = $arr
}

# Copy the formatting (bold/border/center id style + date number format)
# from the last existing data row down onto the two new rows.
$ws.Range("A138:AC138").Copy()
$ws.Range("A139:AC140").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A139").Value2 = 137
$ws.Range("A140").Value2 = 138

Set-RowValues $ws 139 @(
    8117714, "Germany Verbandsliga", "Germany Verbandsliga", 45401.58333333334,
    "SG RotWeiss Thalheim", "Haldensleber SC",
    1, 2, "A",
    2.25, 3.5, 2.6, 2.25, 3.5, 2.6, 0,
    1.775, 2.025, 3, 2, 1.8,
    -1, -1, 1.6, -1, 1.025, 0, 0
)

Set-RowValues $ws 140 @(
    8114461, "Germany Verbandsliga", "Germany Verbandsliga", 45401.625,
    "SG Andernach", "SG 2000 MulheimKarlich",
    2, 3, "A",
    2.35, 4.2, 2.25, 2.625, 4.2, 2, 0.25,
    1.925, 1.875, 4, 2, 1.8,
    -1, -1, 1, -1, 0.875, 1, -1
)
